# Update protocolMap tier values (column C) and mark previously "unknown"
# category rows (column D) as "unrated", per the publish-series update.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("protocolMap")

$tierUpdates = @{
    18 = 1
    22 = 1
    27 = 2
    29 = 2
    30 = 2
    35 = 2
    37 = 2
    40 = 2
    42 = 2
    44 = 2
    45 = 2
    47 = 2
    49 = 2
    50 = 2
    51 = 2
    54 = 2
    55 = 2
    56 = 2
    57 = 2
    58 = 2
    59 = 2
    60 = 2
    61 = 2
    62 = 2
    64 = 3
    65 = 3
    66 = 3
    68 = 3
    71 = 3
    72 = 3
    73 = 3
    74 = 3
    75 = 3
    76 = 3
    77 = 3
    78 = 3
    80 = 3
    81 = 3
    82 = 3
    83 = 3
    84 = 3
    85 = 3
    86 = 3
    87 = 3
    88 = 3
    89 = 3
    90 = 3
    92 = 3
    93 = 3
    96 = 3
    97 = 3
    98 = 3
    99 = 3
    100 = 3
    102 = 3
    103 = 3
    104 = 3
    105 = 3
    106 = 3
    109 = 3
    110 = 3
    111 = 3
    113 = 3
    115 = 3
    117 = 3
    121 = 3
    122 = 3
    124 = 3
    125 = 3
    126 = 3
    128 = 3
    129 = 3
    132 = 3
    133 = 3
    134 = 3
}

foreach ($row in $tierUpdates.Keys) {
    $ws.Cells.Item($row, 3).Value = $tierUpdates[$row]
}

$unratedRows = @(40, 62, 77, 100, 109, 110, 115, 132)
foreach ($row in $unratedRows) {
    $ws.Cells.Item($row, 4).Value = "unrated"
}
